# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
#  1. Insert a new first sheet "Player Info" with the player's basic info.
#  2. Rename MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" / "ODI Bowling"
#     and replace the full scorecard URL values with just the numeric
#     MatchCode that was embedded in the URL's query string.

$wb = $excel.ActiveWorkbook

# Helper: write a value while forcing a text/string cell type (keeps
# numeric-looking strings like match codes or ids from being coerced into
# numbers), then restore the plain "Normal" style so no stray number
# formatting sticks around on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. New "Player Info" sheet, inserted before the existing sheets.
# ---------------------------------------------------------------------
$info = $wb.Worksheets.Add()
$info.Name = "Player Info"

Set-TextValue $info.Range("A1") "ID"
Set-TextValue $info.Range("B1") "NAME"
Set-TextValue $info.Range("C1") "BATTING_HAND"
Set-TextValue $info.Range("D1") "BOWL_STYLE"

$hdr = $info.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

Set-TextValue $info.Range("A2") "6204"
Set-TextValue $info.Range("B2") "Umran Malik"
Set-TextValue $info.Range("C2") "Right Handed"
Set-TextValue $info.Range("D2") "Right Arm Fast"

$info.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK (col D) -> MATCH_CODE
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
# Header text isn't numeric-looking, so a plain assignment keeps the
# existing bold/border header style (s="1") untouched.
$batting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{
    2 = "4669"
    3 = "4673"
    4 = "4676"
    5 = "4682"
    6 = "4685"
    7 = "4687"
    8 = "4689"
    9 = "4697"
}
foreach ($row in $battingCodes.Keys) {
    Set-TextValue $batting.Cells.Item($row, 4) $battingCodes[$row]
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK (col B) -> MATCH_CODE
# ---------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{
    2 = "4669"
    3 = "4676"
    4 = "4682"
    5 = "4685"
    6 = "4687"
    7 = "4689"
    8 = "4697"
}
foreach ($row in $bowlingCodes.Keys) {
    Set-TextValue $bowling.Cells.Item($row, 2) $bowlingCodes[$row]
}
